# ------------------------------------------------------------------
# "2nd testing of object detection using YOLOv8 with 817 data"
#
# Re-layouts the "Testing" sheet:
#   - renames several headers (Percobaan->Testing, Jumlah Dataset->Total
#     Dataset, Jumlah->Total)
#   - inserts a "Total" sub-column inside the Confussion Matrix Testing
#     group (col L), pushing the old "Rasio"/"Akurasi" column into M
#     (relabelled "Accuracy") and the old "Precision" column into N
#   - adds a new "Result Saved" column (O)
#   - fixes the ROUNDUP formulas to round to 0 decimals instead of 2
#   - fills in the 2nd training run (817 images) on row 4
#   - adds a (currently mostly empty) 3rd row of data on row 5
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- undo the old merges that are being reshaped ------------------
$ws.Range("H1:K1").UnMerge()
$ws.Range("L1:L2").UnMerge()
$ws.Range("M1:M2").UnMerge()

# ====================================================================
# Row 1 / Row 2 header text
# ====================================================================
$ws.Range("A1").Value = "Testing"
$ws.Range("B1").Value = "Total Dataset"
$ws.Range("C1").Value = "Train"
$ws.Range("E1").Value = "Test"
$ws.Range("G1").Value = "Epoch"
$ws.Range("H1").Value = "Confussion Matrix Testing"
$ws.Range("M1").Value = "Accuracy"
$ws.Range("N1").Value = "Precision"
$ws.Range("O1").Value = "Result Saved"

$ws.Range("C2").Value = "Rasio"
$ws.Range("D2").Value = "Total"
$ws.Range("E2").Value = "Rasio"
$ws.Range("F2").Value = "Total"
$ws.Range("H2").Value = "TP"
$ws.Range("I2").Value = "FP"
$ws.Range("J2").Value = "FN"
$ws.Range("K2").Value = "TN"
$ws.Range("L2").Value = "Total"

# ====================================================================
# Re-create merges with the new geometry
# ====================================================================
$ws.Range("A1:A2").Merge()
$ws.Range("B1:B2").Merge()
$ws.Range("C1:D1").Merge()
$ws.Range("E1:F1").Merge()
$ws.Range("G1:G2").Merge()
$ws.Range("H1:L1").Merge()
$ws.Range("M1:M2").Merge()
$ws.Range("N1:N2").Merge()
$ws.Range("O1:O2").Merge()

# ====================================================================
# Row 3 - fix the ROUNDUP precision + add the new Total/Result columns
# ====================================================================
$ws.Range("C3").Value = 0.9
$ws.Range("D3").Formula = "=ROUNDUP(C3*B3,0)"
$ws.Range("L3").Formula = "=H3+I3+J3+K3"
$ws.Range("M3").Formula = "=((H3+I3)/(H3+I3+J3+K3))*100"
$ws.Range("N3").Formula = "=(H3/(H3+I3))*100"
$ws.Range("O3").Value = $null

# ====================================================================
# Row 4 - 2nd training run (817 images)
# ====================================================================
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 817
$ws.Range("C4").Value = 0.9
$ws.Range("D4").Formula = "=ROUNDUP(C4*B4,0)"
$ws.Range("E4").Value = 0.1
$ws.Range("F4").Formula = "=B4-D4"
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 45
$ws.Range("I4").Value = 20
$ws.Range("J4").Value = 36
$ws.Range("K4").Value = 0
$ws.Range("L4").Formula = "=H4+I4+J4+K4"
$ws.Range("M4").Formula = "=((H4+I4)/(H4+I4+J4+K4))*100"
$ws.Range("N4").Formula = "=(H4/(H4+I4))*100"
$ws.Range("O4").Value = "runs\detect\train19\weights"

# ====================================================================
# Row 5 - 3rd run placeholder
# ====================================================================
$ws.Range("A5").Value = 3
$ws.Range("C5").Value = 0.9
$ws.Range("D5").Formula = "=ROUNDUP(C5*B5,0)"
$ws.Range("E5").Value = 0.1
$ws.Range("F5").Formula = "=B5-D5"
$ws.Range("L5").Formula = "=H5+I5+J5+K5"
$ws.Range("M5").Formula = "=((H5+I5)/(H5+I5+J5+K5))*100"
$ws.Range("N5").Formula = "=(H5/(H5+I5))*100"

# ====================================================================
# Styling - fill / border / alignment for the header block (rows 1-2)
# ====================================================================
$headerFill = 10216701   # theme green used throughout the sheet's headers (BGR of #FDE49B)

function Set-HeaderBox($addr) {
    $r = $ws.Range($addr)
    $r.Interior.Color = $headerFill
    $r.Font.Name = "Trebuchet MS"
    $r.Font.Size = 10
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
    $r.BorderAround(1, 2)
}

Set-HeaderBox "A1:A2"
Set-HeaderBox "B1:B2"
Set-HeaderBox "C1:D1"
Set-HeaderBox "E1:F1"
Set-HeaderBox "G1:G2"
Set-HeaderBox "H1:L1"
Set-HeaderBox "M1:M2"
Set-HeaderBox "N1:N2"

$rngO = $ws.Range("O1:O2")
$rngO.Interior.Color = $headerFill
$rngO.Font.Name = "Trebuchet MS"
$rngO.Font.Size = 10
$rngO.HorizontalAlignment = -4108
$rngO.VerticalAlignment = -4108
$rngO.WrapText = $true
$rngO.BorderAround(1, 2)

# second sub-header row box (C2:F2, H2:L2) - bordered + filled individually
$rngSub = $ws.Range("C2:F2")
$rngSub.Interior.Color = $headerFill
$rngSub.Font.Name = "Trebuchet MS"
$rngSub.Font.Size = 10
$rngSub.HorizontalAlignment = -4108
$rngSub.VerticalAlignment = -4108
$rngSub.Borders.LineStyle = 1
$rngSub.Borders.Weight = 2

$rngSub2 = $ws.Range("H2:L2")
$rngSub2.Interior.Color = $headerFill
$rngSub2.Font.Name = "Trebuchet MS"
$rngSub2.Font.Size = 10
$rngSub2.HorizontalAlignment = -4108
$rngSub2.VerticalAlignment = -4108
$rngSub2.Borders.Item(7).LineStyle = 1
$rngSub2.Borders.Item(7).Weight = 2
$rngSub2.Borders.Item(10).LineStyle = 1
$rngSub2.Borders.Item(10).Weight = 2
$rngSub2.Borders.Item(9).LineStyle = 1
$rngSub2.Borders.Item(9).Weight = 2
$rngSub2.Borders.Item(8).LineStyle = 0
for ($c = 8; $c -le 11; $c++) {
    $cellRight = $ws.Cells.Item(2, $c)
    $cellRight.Borders.Item(10).LineStyle = 1
    $cellRight.Borders.Item(10).Weight = 2
}

# ====================================================================
# Styling - data rows 3:5 (plain bordered box, no fill)
# ====================================================================
$dataRange = $ws.Range("A3:O5")
$dataRange.Font.Name = "Trebuchet MS"
$dataRange.Font.Size = 10
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108
$dataRange.Borders.LineStyle = 1
$dataRange.Borders.Weight = 2

# M4 carries the 4-decimal numeric format seen on the computed accuracy value
$ws.Range("M4").NumberFormat = "0.0000"

# ====================================================================
# Column widths for the two newly introduced columns
# ====================================================================
$ws.Columns("N").ColumnWidth = 13.29
$ws.Columns("O").ColumnWidth = 28.57

# ====================================================================
# Selection / view state
# ====================================================================
$ws.Range("L15").Select()
